$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.490.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.872.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.81%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.23"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.014"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4798"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.61%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07379"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9400"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.65"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07876"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.875.18"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.450"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.604"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.98"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.017"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008971"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.014"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.94"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.520.31"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.144"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.56"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.14"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08926"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.327"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.219"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.63%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7506"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.710"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02075"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.121"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05305"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5376"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.095"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1523"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.442"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.13%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.61"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.17"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.19"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06112"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9011"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.82%  "
